$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.267595333333334
$ws.Range("H2").Value = 24.802786
$ws.Range("I2").Value = 0.3656664502891759
$ws.Range("J2").Value = 0.3656664502891758
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.07184333333333333
$ws.Range("N2").Value = 0.21553
$ws.Range("O2").Value = 0.04065095086408497
$ws.Range("P2").Value = 0.04065095086408497
$ws.Range("Q2").Value = 0.5939716073977778
$ws.Range("R2").Value = 5.34574446658
$ws.Range("S2").Value = 0.01486468890334966
$ws.Range("T2").Value = 0.01486468890334965
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.267595333333334
$ws.Range("H3").Value = 24.802786
$ws.Range("I3").Value = 0.3656664502891759
$ws.Range("J3").Value = 0.3656664502891758
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.315200666666667
$ws.Range("N3").Value = 3.945602
$ws.Range("O3").Value = 0.7441770195853729
$ws.Range("P3").Value = 0.7441770195853727
$ws.Range("Q3").Value = 10.87354689413022
$ws.Range("R3").Value = 97.861922047172
$ws.Range("S3").Value = 0.2721205691385618
$ws.Range("T3").Value = 0.2721205691385618
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.267595333333334
$ws.Range("H4").Value = 24.802786
$ws.Range("I4").Value = 0.3656664502891759
$ws.Range("J4").Value = 0.3656664502891758
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3802783333333333
$ws.Range("N4").Value = 1.140835
$ws.Range("O4").Value = 0.2151720295505423
$ws.Range("P4").Value = 0.2151720295505423
$ws.Range("Q4").Value = 3.143987374034445
$ws.Range("R4").Value = 28.29588636631
$ws.Range("S4").Value = 0.07868119224726446
$ws.Range("T4").Value = 0.07868119224726444
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.10830433333333
$ws.Range("H5").Value = 33.324913
$ws.Range("I5").Value = 0.491307816908375
$ws.Range("J5").Value = 0.491307816908375
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.07184333333333333
$ws.Range("N5").Value = 0.21553
$ws.Range("O5").Value = 0.04065095086408497
$ws.Range("P5").Value = 0.04065095086408497
$ws.Range("Q5").Value = 0.7980576109877776
$ws.Range("R5").Value = 7.182518498889999
$ws.Range("S5").Value = 0.01997212992428321
$ws.Range("T5").Value = 0.01997212992428321
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.10830433333333
$ws.Range("H6").Value = 33.324913
$ws.Range("I6").Value = 0.491307816908375
$ws.Range("J6").Value = 0.491307816908375
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.315200666666667
$ws.Range("N6").Value = 3.945602
$ws.Range("O6").Value = 0.7441770195853729
$ws.Range("P6").Value = 0.7441770195853727
$ws.Range("Q6").Value = 14.60964926473622
$ws.Range("R6").Value = 131.486843382626
$ws.Range("S6").Value = 0.3656199868858706
$ws.Range("T6").Value = 0.3656199868858705
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.10830433333333
$ws.Range("H7").Value = 33.324913
$ws.Range("I7").Value = 0.491307816908375
$ws.Range("J7").Value = 0.491307816908375
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3802783333333333
$ws.Range("N7").Value = 1.140835
$ws.Range("O7").Value = 0.2151720295505423
$ws.Range("P7").Value = 0.2151720295505423
$ws.Range("Q7").Value = 4.224247458039444
$ws.Range("R7").Value = 38.01822712235499
$ws.Range("S7").Value = 0.1057157000982213
$ws.Range("T7").Value = 0.1057157000982213
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.233763666666666
$ws.Range("H8").Value = 9.701290999999999
$ws.Range("I8").Value = 0.1430257328024492
$ws.Range("J8").Value = 0.1430257328024492
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.07184333333333333
$ws.Range("N8").Value = 0.21553
$ws.Range("O8").Value = 0.04065095086408497
$ws.Range("P8").Value = 0.04065095086408497
$ws.Range("Q8").Value = 0.2323243610255555
$ws.Range("R8").Value = 2.09091924923
$ws.Range("S8").Value = 0.005814132036452109
$ws.Range("T8").Value = 0.005814132036452109
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.233763666666666
$ws.Range("H9").Value = 9.701290999999999
$ws.Range("I9").Value = 0.1430257328024492
$ws.Range("J9").Value = 0.1430257328024492
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.315200666666667
$ws.Range("N9").Value = 3.945602
$ws.Range("O9").Value = 0.7441770195853729
$ws.Range("P9").Value = 0.7441770195853727
$ws.Range("Q9").Value = 4.253048130242444
$ws.Range("R9").Value = 38.277433172182
$ws.Range("S9").Value = 0.1064364635609406
$ws.Range("T9").Value = 0.1064364635609406
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.233763666666666
$ws.Range("H10").Value = 9.701290999999999
$ws.Range("I10").Value = 0.1430257328024492
$ws.Range("J10").Value = 0.1430257328024492
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3802783333333333
$ws.Range("N10").Value = 1.140835
$ws.Range("O10").Value = 0.2151720295505423
$ws.Range("P10").Value = 0.2151720295505423
$ws.Range("Q10").Value = 3.143987374034445
$ws.Range("R10").Value = 11.067572317985
$ws.Range("S10").Value = 0.03077513720505657
$ws.Range("T10").Value = 0.03077513720505657
